$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- TC3 continuation: patient registration flow ---
$ws.Range("A19").Value = "1. Click the patient search box (Enter name, phone or ID)"
$ws.Range("A20").Value = "2. Click REGISTER AS A NEW PATIENT button"
$ws.Range("A21").Value = "3. Fillout all the fields"
$ws.Range("A22").Value = "4. Click Next button"
$ws.Range("A23").Value = "5. Select NO for all disease"
$ws.Range("A24").Value = "6. Select YES for all disease"
$ws.Range("A25").Value = "7. Click save"
$ws.Range("A26").Value = "8. Click Not Now"
$ws.Range("A27").Value = "9. Validate patient if exist via search"

# --- TC4 header (highlighted like TC1/TC2/TC3) ---
$ws.Range("A28").Value = "TC4"
$ws.Range("A28").Interior.Color = 65535

# --- TC4 steps: medicine selection flow ---
$ws.Range("A29").Value = "1. Click Medicine (+) button"
$ws.Range("A30").Value = "2. Select multiple kinds of existing medicines and their dosage"
$ws.Range("A31").Value = "3. Click Save"
$ws.Range("A32").Value = "4. Validate if the selected medicines reflected in the medicine section"

# --- TC5 header (highlighted like TC1/TC2/TC3/TC4) ---
$ws.Range("A33").Value = "TC5"
$ws.Range("A33").Interior.Color = 65535

# Reflect the final cursor/selection position left by the author (row 32 area, column H)
[void]$ws.Range("H32").Select()
